$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.06398557126522064, 0.982882022857666, 0.0196644626557827, 0.996960461139679),
    @(0.01011102832853794, 0.9984486103057861, 0.0102210994809866, 0.9976118206977844),
    @(0.004525858443230391, 0.9988717436790466, 0.00237832753919065, 0.9989144802093506),
    @(0.002398055279627442, 0.9994182586669922, 0.001493193209171295, 0.9991315603256226),
    @(0.002402569865807891, 0.9995239973068237, 0.0004406488442327827, 1),
    @(0.001922777388244867, 0.9995063543319702, 0.0001149797026300803, 1),
    @(0.001402237918227911, 0.9996474385261536, 0.0001106618583435193, 1),
    @(0.001023262972012162, 0.9997355341911316, 0.0000718216979294084, 1),
    @(0.0009293231996707618, 0.9997884631156921, 0.0005154672544449568, 0.999782919883728),
    @(0.0004536720807664096, 0.9998942017555237, 0.00004631537740351632, 1),
    @(0.000903626496437937, 0.9998236894607544, 0.00007343941979343072, 1),
    @(0.001016956521198153, 0.9997884631156921, 0.00004404792707646266, 1),
    @(0.001087073935195804, 0.9997884631156921, 0.0002127236075466499, 1),
    @(0.0003716057108249515, 0.9999294877052307, 0.00002430202403047588, 1),
    @(0.0004443767538759857, 0.9999118447303772, 0.002566620940342546, 0.9995657801628113),
    @(0.0004621722910087556, 0.9998766183853149, 0.0000008838957796797331, 1),
    @(0.0001379724126309156, 0.999964714050293, 0.0000006284655569288589, 1),
    @(0.0007585492567159235, 0.9998413324356079, 0.00001520953719591489, 1),
    @(0.0004452612774912268, 0.9999118447303772, 0.0002602652239147574, 0.999782919883728),
    @(0.0005967464530840516, 0.9998942017555237, 0.0000005659768476107274, 1),
    @(0.0007819700986146927, 0.9998589754104614, 0.00005221238825470209, 1),
    @(0.0001364952477160841, 0.9999471306800842, 0.000112151654320769, 1),
    @(0.0003777554084081203, 0.9999294877052307, 0.0000006930375775482389, 1),
    @(0.001085414434783161, 0.9997708201408386, 0.0001733792305458337, 1),
    @(0.000124375888844952, 0.9999471306800842, 0.000001046317720465595, 1),
    @(0.00006842263246653602, 0.9999823570251465, 0.0003234297910239547, 0.999782919883728),
    @(0.001097306725569069, 0.9998413324356079, 0.000003128557864329196, 1),
    @(0.0002468824677634984, 0.9999294877052307, 0.0000002649328791903827, 1),
    @(0.0001321008458035067, 0.9999471306800842, 0.000000352391111846373, 1),
    @(0.0007811344694346189, 0.9998236894607544, 0.000000314151265001783, 1),
    @(0.0005374005413614213, 0.9998589754104614, 0.000003507638666633284, 1),
    @(0.0001660581328906119, 0.9999471306800842, 0.000001364840045425808, 1),
    @(0.0001811858237488195, 0.9999823570251465, 0.0000001127746998008661, 1),
    @(0.00005568853157456033, 0.999964714050293, 0.0000002738958073678077, 1),
    @(0.00004675958189181983, 0.999964714050293, 0.000000002743415050687759, 1),
    @(0.001124259433709085, 0.9998942017555237, 0.0001055254251696169, 1),
    @(0.0001433350407751277, 0.9999294877052307, 0.0000005002899570172303, 1),
    @(0.0002093774382956326, 0.999964714050293, 0.000000009576041293257731, 1),
    @(0.0002415008493699133, 0.9999294877052307, 0.0000004940908411299461, 1),
    @(0.00004266534597263671, 0.9999823570251465, 0.0000001291627995669842, 1),
    @(0.0001516610936960205, 0.999964714050293, 0.00000008166973231027441, 1),
    @(0.0003442777961026877, 0.9999118447303772, 0.000002033769987974665, 1),
    @(0.0008760616765357554, 0.9998766183853149, 0.0000006563233796441637, 1),
    @(0.00007451117562595755, 0.999964714050293, 0.000000008411382701467573, 1),
    @(0.0003570486733224243, 0.9999823570251465, 0.0000001483272455971019, 1),
    @(0.0006082363543100655, 0.9998942017555237, 0.0001057337940437719, 1),
    @(0.0003655508917290717, 0.9998766183853149, 0.00000005636593769509091, 1),
    @(0.00004778323273058049, 0.9999823570251465, 0.0000004603427896654466, 1),
    @(0.00003219924838049337, 1, 0.00000001164643315121339, 1),
    @(0.0000194378608284751, 1, 0.00000002823573197474616, 1),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
